$d = $word.ActiveDocument

# Locate the paragraph that ends the "Key Assignment" section, i.e. the one
# mentioning switching/moving to a desktop by number. The new content is
# inserted right after it, before the "Default Browser Activation" heading.
$anchorText = "You can also switch/move to a particular desktop by number."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$anchorText*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find anchor paragraph for insertion."
}

# Insert a new empty paragraph right after the anchor paragraph, then fill it
# in as the new "Window Activation on Switch From Empty Desktop" heading.
$target.Range.InsertParagraphAfter()
$headingPara = $target.Next()
$headingPara.Range.Text = "Window Activation on Switch From Empty Desktop"
$headingPara.Style = "Heading 2"

# Insert another new paragraph after the heading for the descriptive body text.
$headingPara.Range.InsertParagraphAfter()
$bodyPara = $headingPara.Next()
$bodyPara.Style = "Normal"
$bodyPara.Range.Text = "Virtual Desktop Grid Switcher fixes an issue in Windows 10 where switching from a desktop which is empty to a desktop which had an activate window does not reactivate that window."

Write-Output "Inserted 'Window Activation on Switch From Empty Desktop' section."
